$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.588.00"
$ws.Range("E2").Value = "  -3.17%  "
$ws.Range("D3").Value = "'1.848.55"
$ws.Range("E3").Value = "  -3.91%  "
$ws.Range("E4").Value = "  -1.05%  "
$ws.Range("D5").Value = "'335.58"
$ws.Range("E5").Value = "  +2.87%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("D7").Value = "'0.4662"
$ws.Range("E7").Value = "  -3.07%  "
$ws.Range("D8").Value = "'0.3898"
$ws.Range("E8").Value = "  -3.65%  "
$ws.Range("D9").Value = "'46.14"
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").Value = "'0.07905"
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("D11").Value = "'0.9786"
$ws.Range("E11").Value = "  -2.92%  "
$ws.Range("D12").Value = "'22.26"
$ws.Range("E12").Value = "  -6.52%  "
$ws.Range("D13").Value = "'1.851.19"
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("D14").Value = "'5.820"
$ws.Range("E14").Value = "  -4.38%  "
$ws.Range("D15").Value = "'6.984"
$ws.Range("E15").Value = "  -4.47%  "
$ws.Range("D16").Value = "'0.06911"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "'87.53"
$ws.Range("E18").Value = "  -4.37%  "
$ws.Range("D19").Value = "'0.00001002"
$ws.Range("E19").Value = "  -3.54%  "
$ws.Range("D20").Value = "'17.06"
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").Value = "'28.594.14"
$ws.Range("E22").Value = "  -3.14%  "
$ws.Range("D23").Value = "'5.387"
$ws.Range("E23").Value = "  -4.85%  "
$ws.Range("D24").Value = "'11.26"
$ws.Range("E24").Value = "  -6.02%  "
$ws.Range("D25").Value = "'2.161"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").Value = "'2.053.61"
$ws.Range("E26").Value = "  -4.40%  "
$ws.Range("D27").Value = "'153.22"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("D28").Value = "'19.42"
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("E29").Value = "  -5.17%  "
$ws.Range("D30").Value = "'2.023"
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("D31").Value = "'117.34"
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("D32").Value = "'0.9692"
$ws.Range("E32").Value = "  -4.23%  "
$ws.Range("D33").Value = "'0.09333"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").Value = "'5.352"
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("D35").Value = "'3.482"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").Value = "'1.342"
$ws.Range("E36").Value = "  -3.02%  "
$ws.Range("D37").Value = "'0.06154"
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("D38").Value = "'0.02197"
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("D39").Value = "'1.169"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5703"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.672"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("D42").Value = "'10.13"
$ws.Range("E42").Value = "  -5.55%  "
$ws.Range("D43").Value = "'0.1790"
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").Value = "'2.412"
$ws.Range("E44").Value = "  -2.55%  "
$ws.Range("D45").Value = "'1.247"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5378"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'11.74"
$ws.Range("E47").Value = "  -5.32%  "
$ws.Range("D48").Value = "'0.07096"
$ws.Range("E48").Value = "  -4.98%  "
$ws.Range("D49").Value = "'1.902"
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("D50").Value = "'113.35"
$ws.Range("E50").Value = "  -3.71%  "
$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  -1.14%  "
